$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "Main"
$ws.Name = "Main"

# Turn off gridlines for the view
$excel.ActiveWindow.DisplayGridlines = $false

# ---------------- Values ----------------
# Column B labels, top to bottom
$ws.Range("B2").Value = "Ticker"
$ws.Range("B3").Value = "Price"
$ws.Range("B4").Value = "S/O"
$ws.Range("B5").Value = "Mkt Cap"
$ws.Range("B6").Value = "Cash"
$ws.Range("B7").Value = "Debt"
$ws.Range("B8").Value = "EV"

# Column C values/formulas
$ws.Range("C2").Value = "ABVX"
$ws.Range("C3").Value = 70.5
$ws.Range("C4").Value = 77400000
$ws.Range("C5").Formula = "=C4*C3"
$ws.Range("C6").Value = 170671000
$ws.Range("C7").Formula = "=(70645+93999)*1000"
$ws.Range("C8").Formula = "=C5+C7-C6"

# "In Euros" note
$ws.Range("E2").Value = "In Euros"

# Pipeline summary (columns G and I)
$ws.Range("G2").Value = "Obefazimod "
$ws.Range("G3").Value = "ABTECT - PIII"
$ws.Range("G4").Value = "UC"
$ws.Range("I2").Value = "Obefazimod "
$ws.Range("I3").Value = "ENHANCE-CD - PIIb"
$ws.Range("I4").Value = "Chrons"

# ---------------- Fonts ----------------
# Base font "Aptos Serif" applied cell-by-cell so no empty filler
# cells get created between the populated columns.
$dataCells = @("B2","B3","B4","B5","B6","B7","B8",
               "C2","C3","C4","C5","C6","C7","C8",
               "E2",
               "G2","G3","G4",
               "I2","I3","I4")
foreach ($addr in $dataCells) {
    $ws.Range($addr).Font.Name = "Aptos Serif"
}

# "In Euros" header is bold + underlined
$ws.Range("E2").Font.Bold = $true
$ws.Range("E2").Font.Underline = $true

# ---------------- Number formats ----------------
$ws.Range("C4:C8").NumberFormat = "#,##0"

# ---------------- Alignment ----------------
$ws.Range("C2").HorizontalAlignment = -4152

# ---------------- Column widths ----------------
$ws.Columns("A:B").ColumnWidth = 9.140625
$ws.Columns("C:C").ColumnWidth = 14.28515625

# ---------------- Selection ----------------
$ws.Range("I5").Select() | Out-Null
